# COVID-19 TW PlaceCode workbook update
# Adds three new daily-report date columns (L, M, N) on row 1, fills in the
# matching "case count" header cells on row 2, populates the new J-column
# (2020-04-04) case counts for every place row, and removes the now-unused
# placeholder cells in columns M/N that used to reserve space for this data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeBook")

# --- Row 1: extend the date header with three more days -------------------
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$ws.Range("L1").Value = 43927
$ws.Range("M1").Value = 43928
$ws.Range("N1").Value = 43929

# --- Row 2: extend the "病例數" (case count) header label ------------------
$ws.Range("K2").Copy()
$ws.Range("L2:N2").PasteSpecial(-4122)
$ws.Range("L2").Value = "病例數"
$ws.Range("M2").Value = "病例數"
$ws.Range("N2").Value = "病例數"

# --- Column J (2020-04-04 counts): copy number style from column I, -------
# --- then fill in each place's case count for that day --------------------
$ws.Range("I3").Copy()
$ws.Range("J3:J24").PasteSpecial(-4122)

$jValues = @{
    3  = 104
    4  = 34
    5  = 13
    6  = 27
    7  = 3
    8  = 9
    9  = 3
    10 = 83
    11 = 36
    12 = 6
    13 = 2
    14 = 3
    15 = 17
    16 = 2
    17 = 4
    18 = 2
    19 = 7
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
}

foreach ($r in $jValues.Keys) {
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}

# --- Remove the now-obsolete placeholder cells in columns M/N -------------
for ($r = 3; $r -le 24; $r++) {
    $ws.Cells.Item($r, 13).Clear()
}
$ws.Cells.Item(25, 13).Clear()
$ws.Cells.Item(25, 14).Clear()

# --- Row 25 (totals) shrinks back to its natural height now that the bold
# --- placeholder cells are gone -------------------------------------------
$ws.Rows.Item(25).RowHeight = 17

# --- Update the remembered selection to match the newly-entered range -----
[void]$ws.Range("K2:N2").Select()

Write-Output "edit complete"
